$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay text so numeric-looking price strings are not
# reinterpreted as numbers (losing exact formatting / precision).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.325.07"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "1.875.37"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "0.7123"
$ws.Range("E5").Value = "  -0.21%  "

# Row 6
$ws.Range("D6").Value = "241.99"
$ws.Range("E6").Value = "  +0.63%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "0.3100"
$ws.Range("E8").Value = "  +0.81%  "

# Row 9
$ws.Range("D9").Value = "0.07731"
$ws.Range("E9").Value = "  -0.31%  "

# Row 10
$ws.Range("D10").Value = "24.84"
$ws.Range("E10").Value = "  -1.04%  "

# Row 11
$ws.Range("D11").Value = "0.08529"
$ws.Range("E11").Value = "  +3.26%  "

# Row 12
$ws.Range("D12").Value = "1.877.00"
$ws.Range("E12").Value = "  -0.34%  "

# Row 13
$ws.Range("E13").Value = "  -0.30%  "

# Row 14
$ws.Range("D14").Value = "0.7097"

# Row 15
$ws.Range("E15").Value = "  +1.16%  "

# Row 16
$ws.Range("D16").Value = "29.318.59"
$ws.Range("E16").Value = "  +0.37%  "

# Row 17
$ws.Range("D17").Value = "0.000008211"
$ws.Range("E17").Value = "  +5.30%  "

# Row 18
$ws.Range("D18").Value = "6.008"
$ws.Range("E18").Value = "  +2.42%  "

# Row 19
$ws.Range("D19").Value = "241.55"
$ws.Range("E19").Value = "  -1.29%  "

# Row 20
$ws.Range("D20").Value = "2.134.11"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  +0.79%  "

# Row 22
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").Value = "7.805"
$ws.Range("E23").Value = "  -2.43%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").Value = "0.1600"
$ws.Range("E25").Value = "  +0.63%  "

# Row 26
$ws.Range("D26").Value = "163.39"
$ws.Range("E26").Value = "  +0.81%  "

# Row 27
$ws.Range("D27").Value = "9.030"
$ws.Range("E27").Value = "  +1.39%  "

# Row 28
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").Value = "  +0.63%  "

# Row 29
$ws.Range("D29").Value = "1.514"
$ws.Range("E29").Value = "  +1.25%  "

# Row 30
$ws.Range("D30").Value = "4.397"
$ws.Range("E30").Value = "  -0.49%  "

# Row 31
$ws.Range("D31").Value = "4.313"
$ws.Range("E31").Value = "  +2.22%  "

# Row 32
$ws.Range("D32").Value = "1.285"
$ws.Range("E32").Value = "  -2.65%  "

# Row 33
$ws.Range("D33").Value = "0.05267"
$ws.Range("E33").Value = "  +1.63%  "

# Row 34
$ws.Range("D34").Value = "1.933"
$ws.Range("E34").Value = "  +1.07%  "

# Row 35
$ws.Range("D35").Value = "1.174"

# Row 36
$ws.Range("D36").Value = "0.7452"
$ws.Range("E36").Value = "  +2.45%  "

# Row 37
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +0.40%  "

# Row 38
$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +0.74%  "

# Row 39
$ws.Range("D39").Value = "2.720"
$ws.Range("E39").Value = "  +1.30%  "

# Row 40
$ws.Range("D40").Value = "1.180.31"

# Row 41
$ws.Range("D41").Value = "6.384"
$ws.Range("E41").Value = "  +3.68%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "72.85"
$ws.Range("E42").Value = "  +0.63%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8863"
$ws.Range("E43").Value = "  -2.17%  "

# Row 44
$ws.Range("D44").Value = "106.43"
$ws.Range("E44").Value = "  +4.69%  "

# Row 45
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.05%  "

# Row 46
$ws.Range("D46").Value = "2.030.10"
$ws.Range("E46").Value = "  +1.31%  "

# Row 47
$ws.Range("D47").Value = "1.812"
$ws.Range("E47").Value = "  +2.65%  "

# Row 48
$ws.Range("D48").Value = "0.5208"
$ws.Range("E48").Value = "  -0.09%  "

# Row 49
$ws.Range("D49").Value = "0.00000000122"
$ws.Range("E49").Value = "  +1.37%  "

# Row 50
$ws.Range("D50").Value = "9.390"
$ws.Range("E50").Value = "  +0.96%  "

# Row 51
$ws.Range("D51").Value = "0.4315"
$ws.Range("E51").Value = "  +1.30%  "
